$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# ---------------------------------------------------------------------
# 1) Merge "User " + "Behaviour" + " Data:" into a single run's text.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("User Behaviour Data:", $true, $false, $false, $false, $false, $true, `
    $wdFindContinue, $false, "User Behaviour Data:", $wdReplaceOne)
if (-not $found) { throw "Could not find 'User Behaviour Data:'" }

# ---------------------------------------------------------------------
# 2) Merge ". These models " + "analyse" + " the vast dataset..." into one run.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(". These models analyse the vast dataset mentioned above to make personalized product recommendations.", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, `
    ". These models analyse the vast dataset mentioned above to make personalized product recommendations.", $wdReplaceOne)
if (-not $found) { throw "Could not find 'These models analyse' sentence" }

# ---------------------------------------------------------------------
# 3) "Amazon uses your buying history..." -> underline the middle clause.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("buying history, search history, time spent searching for a product, and many other factors to provide personalized product recommendations ", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) { throw "Could not find buying history clause" }
$rng.Font.Underline = 1

# ---------------------------------------------------------------------
# 4) "Amazon uses your data to identify..." -> underline the middle clause.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("uses your data to identify your interests and show ads about the product and services based on your interests", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) { throw "Could not find 'uses your data to identify' clause" }
$rng.Font.Underline = 1

# ---------------------------------------------------------------------
# 5) "Amazon identifies the buying habits of each customer..." -> underline two clauses.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("identifies the buying habits of each customer", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) { throw "Could not find 'identifies the buying habits' clause" }
$rng.Font.Underline = 1

$rng = $d.Content
$found = $rng.Find.Execute("helps identify the locations with the most frequent buyers and employ more delivery partners in such areas", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) { throw "Could not find 'helps identify the locations' clause" }
$rng.Font.Underline = 1

# ---------------------------------------------------------------------
# 6) Alexa paragraph -> underline three clauses.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("quick and accurate answers to customer questions", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) { throw "Could not find 'quick and accurate answers' clause" }
$rng.Font.Underline = 1

$rng = $d.Content
$found = $rng.Find.Execute("questions about products, place orders, and provide information about the weather, traffic, and more", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) { throw "Could not find 'questions about products' clause" }
$rng.Font.Underline = 1

$rng = $d.Content
$found = $rng.Find.Execute("able to reduce the burden on its human customer service representatives and provide a more efficient and convenient experience for customers.", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
if (-not $found) { throw "Could not find 'able to reduce the burden' clause" }
$rng.Font.Underline = 1

# ---------------------------------------------------------------------
# 7) Merge "By leveraging user " + "behaviour" + " data, product..." into one run.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("By leveraging user behaviour data, product information, and advanced machine learning models, Amazon can suggest products tailored to individual preferences, ultimately driving customer satisfaction, increasing sales, and solidifying its position as a global e-commerce giant.", `
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, `
    "By leveraging user behaviour data, product information, and advanced machine learning models, Amazon can suggest products tailored to individual preferences, ultimately driving customer satisfaction, increasing sales, and solidifying its position as a global e-commerce giant.", $wdReplaceOne)
if (-not $found) { throw "Could not find 'By leveraging user behaviour data' sentence" }

Write-Host "Done."
